$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "test_book_2"
$ws.Range("E3").Value = "test_location_2"
$ws.Range("F3").Value = "test_subsidiary_2"
$ws.Range("G3").Value = "test_organisation_2"
$ws.Range("H3").Value = "test_risk_class_2"
$ws.Range("I3").Value = "test_var_hierarchy_2"
$ws.Range("J3").Value = "var_facet_2"
$ws.Range("K3").Value = "test_asset_allocation_2"

$ws.Range("K3").Select()
